$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 1058:1059, pushing the existing data
# (previously rows 1058-1157) down to 1060-1159. This also duplicates
# the trailing rows correctly since the sheet's used range grows by 2.
$ws.Rows("1058:1059").Insert()

# Populate the new row 1058 (Primera, new weekly observation dated 45166)
$ws.Range("A1058").Value = 3
$ws.Range("B1058").Value = "Femacal de La Calera"
$ws.Range("C1058").Value = "Coquimbo"
$ws.Range("D1058").Value = 45166
$ws.Range("E1058").Value = 5
$ws.Range("F1058").Value = 100112008
$ws.Range("G1058").Value = "Coliflor"
$ws.Range("H1058").Value = "Sin especificar"
$ws.Range("I1058").Value = "Primera"
$ws.Range("J1058").Value = 2500
$ws.Range("K1058").Value = 650
$ws.Range("L1058").Value = 700
$ws.Range("M1058").Value = 676
$ws.Range("N1058").Value = "$/unidad"
$ws.Range("O1058").Value = "Provincia de Quillota"
$ws.Range("P1058").Value = 676
$ws.Range("Q1058").Value = 1
$ws.Range("R1058").Value = "Hortaliza"

# Populate the new row 1059 (Segunda, new weekly observation dated 45166)
$ws.Range("A1059").Value = 3
$ws.Range("B1059").Value = "Femacal de La Calera"
$ws.Range("C1059").Value = "Coquimbo"
$ws.Range("D1059").Value = 45166
$ws.Range("E1059").Value = 5
$ws.Range("F1059").Value = 100112008
$ws.Range("G1059").Value = "Coliflor"
$ws.Range("H1059").Value = "Sin especificar"
$ws.Range("I1059").Value = "Segunda"
$ws.Range("J1059").Value = 1100
$ws.Range("K1059").Value = 550
$ws.Range("L1059").Value = 550
$ws.Range("M1059").Value = 550
$ws.Range("N1059").Value = "$/unidad"
$ws.Range("O1059").Value = "Provincia de Quillota"
$ws.Range("P1059").Value = 550
$ws.Range("Q1059").Value = 1
$ws.Range("R1059").Value = "Hortaliza"
